# Applies the "Test with a variable production cost" edit:
#  - Productdata sheet: add a new "VariableCosts" column (J) with a header
#    and per-product values, and update the StandardDevDemands value for
#    the first product (H2).
#  - ForecastedAverageDemand sheet: update several forecasted demand values
#    (column B) that shifted because of the new random draws used for the
#    variable cost column.
#  - ForcastedStandardDeviation sheet: update the corresponding standard
#    deviation values (column B) for the same reason.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Productdata sheet: add VariableCosts column
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Productdata")

# Copy the style of the last existing header cell (I1, bold/bordered) onto
# the new header cell J1, then set its text.
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 10))
$ws.Cells.Item(1, 10).Value = "VariableCosts"

# Updated StandardDevDemands value for the first product.
$ws.Cells.Item(2, 8).Value = 21.00793494348941

# New VariableCosts values, one per product row (2-11).
$ws.Cells.Item(2, 10).Value = 9
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(5, 10).Value = 2
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(11, 10).Value = 0

# ---------------------------------------------------------------------
# ForecastedAverageDemand sheet: refresh sampled demand values
# ---------------------------------------------------------------------
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")

$wsAvg.Cells.Item(2, 2).Value = 0
$wsAvg.Cells.Item(3, 2).Value = 0
$wsAvg.Cells.Item(4, 2).Value = 0
$wsAvg.Cells.Item(5, 2).Value = 90
$wsAvg.Cells.Item(6, 2).Value = 101
$wsAvg.Cells.Item(7, 2).Value = 87
$wsAvg.Cells.Item(8, 2).Value = 107
$wsAvg.Cells.Item(9, 2).Value = 105
$wsAvg.Cells.Item(10, 2).Value = 97
$wsAvg.Cells.Item(11, 2).Value = 104
$wsAvg.Cells.Item(12, 2).Value = 116
$wsAvg.Cells.Item(13, 2).Value = 105
$wsAvg.Cells.Item(14, 2).Value = 93
$wsAvg.Cells.Item(15, 2).Value = 98
$wsAvg.Cells.Item(16, 2).Value = 107
$wsAvg.Cells.Item(17, 2).Value = 95
$wsAvg.Cells.Item(19, 2).Value = 108
$wsAvg.Cells.Item(20, 2).Value = 108
$wsAvg.Cells.Item(21, 2).Value = 104
$wsAvg.Cells.Item(22, 2).Value = 100
$wsAvg.Cells.Item(23, 2).Value = 83

# ---------------------------------------------------------------------
# ForcastedStandardDeviation sheet: refresh sampled std-dev values
# ---------------------------------------------------------------------
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")

$wsStd.Cells.Item(2, 2).Value = -0
$wsStd.Cells.Item(3, 2).Value = -0
$wsStd.Cells.Item(5, 2).Value = 11.25
$wsStd.Cells.Item(6, 2).Value = 18.9375
$wsStd.Cells.Item(7, 2).Value = 19.03125
$wsStd.Cells.Item(8, 2).Value = 25.078125
$wsStd.Cells.Item(9, 2).Value = 25.4296875
$wsStd.Cells.Item(10, 2).Value = 23.87109375
$wsStd.Cells.Item(11, 2).Value = 25.796875
$wsStd.Cells.Item(12, 2).Value = 28.88671875
$wsStd.Cells.Item(13, 2).Value = 26.19873046875
$wsStd.Cells.Item(14, 2).Value = 23.227294921875
$wsStd.Cells.Item(15, 2).Value = 24.488037109375
$wsStd.Cells.Item(16, 2).Value = 26.74346923828125
$wsStd.Cells.Item(17, 2).Value = 23.74710083007812
$wsStd.Cells.Item(19, 2).Value = 26.99917602539062
$wsStd.Cells.Item(20, 2).Value = 26.99958801269531
$wsStd.Cells.Item(21, 2).Value = 25.99980163574219
$wsStd.Cells.Item(22, 2).Value = 24.99990463256836
$wsStd.Cells.Item(23, 2).Value = 20.74996042251587
